$p = $ppt.ActivePresentation

# --- Slide 9: remove the empty "Content Placeholder 2" box, keep the Title ---
$s9 = $p.Slides.Item(9)
for ($i = $s9.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s9.Shapes.Item($i)
    if ($sh.Name -eq "Content Placeholder 2") {
        $sh.Delete()
    }
}

# --- Slide 14: remove the empty "Title 1" box, keep the picture ---
$s14 = $p.Slides.Item(14)
for ($i = $s14.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s14.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") {
        $sh.Delete()
    }
}
